$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "X/Y Percentage Error" columns (D, E) hold numeric-looking text.
# Force text formatting first so Excel keeps them as strings instead of
# coercing to numbers.
$ws.Range("D2:E11").NumberFormat = "@"

# Row 2
$ws.Range("C2").Value = "(32536.95729,32868.77145)"
$ws.Range("D2").Value = "-4067044.66085"
$ws.Range("E2").Value = "-6573694.2909"
$ws.Range("F2").Value = "Fail"

# Row 3
$ws.Range("B3").Value = "(0.3082,0.27783)"
$ws.Range("C3").Value = "(0.30945,0.27272)"
$ws.Range("D3").Value = "-0.1559"
$ws.Range("E3").Value = "1.02215"
$ws.Range("F3").Value = "Pass"

# Row 4
$ws.Range("B4").Value = "(0.05246,0.01999)"
$ws.Range("C4").Value = "(0.05343,-2e-05)"
$ws.Range("D4").Value = "-0.12176"
$ws.Range("E4").Value = "4.00208"
$ws.Range("F4").Value = "Pass"

# Row 5 (Pass/Fail column F unchanged)
$ws.Range("B5").Value = "(0.39949,0.19011)"
$ws.Range("C5").Value = "(0.39937,-0.06933)"
$ws.Range("D5").Value = "0.01468"
$ws.Range("E5").Value = "51.8881"

# Row 6
$ws.Range("B6").Value = "(0.78175,0.31904)"
$ws.Range("C6").Value = "(2113877597.61496,1424524248.54673)"
$ws.Range("D6").Value = "-264234699604.151"
$ws.Range("E6").Value = "-284904849645.5389"
$ws.Range("F6").Value = "Fail"

# Row 7 (Pass/Fail column F unchanged)
$ws.Range("B7").Value = "(0.2981,0.47984)"
$ws.Range("C7").Value = "(0.29489,0.47671)"
$ws.Range("D7").Value = "0.40066"
$ws.Range("E7").Value = "0.62517"

# Row 8
$ws.Range("B8").Value = "(0.50199,0.23276)"
$ws.Range("C8").Value = "(52009.71262,63261.82058)"
$ws.Range("D8").Value = "-6501151.32836"
$ws.Range("E8").Value = "-12652317.56278"
$ws.Range("F8").Value = "Fail"

# Row 9
$ws.Range("B9").Value = "(0.52727,0.4496)"
$ws.Range("C9").Value = "(1.08143,1.41821)"
$ws.Range("D9").Value = "-69.26966"
$ws.Range("E9").Value = "-193.72171"
$ws.Range("F9").Value = "Fail"

# Row 10
$ws.Range("B10").Value = "(0.72516,0.45098)"
$ws.Range("C10").Value = "(138703.61521,138089.69658)"
$ws.Range("D10").Value = "-17337861.25595"
$ws.Range("E10").Value = "-27617849.11903"
$ws.Range("F10").Value = "Fail"

# Row 11
$ws.Range("B11").Value = "(0.56015,0.14296)"
$ws.Range("C11").Value = "(30325065.49979,20431960.68204)"
$ws.Range("D11").Value = "-3790633117.45519"
$ws.Range("E11").Value = "-4086392107.81594"
$ws.Range("F11").Value = "Fail"
